# Insert a new data row at row 7 (pushing the existing rows 7-52 down to 8-53)
# and populate it with the new "Inferno" record dated 2021-11-03.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(7).Insert()

$ws.Cells.Item(7, 1).Value = 11
$ws.Cells.Item(7, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(7, 3).Value = "Bíobío"
$ws.Cells.Item(7, 4).Value = 44503
$ws.Cells.Item(7, 5).Value = 8
$ws.Cells.Item(7, 6).Value = 100112021
$ws.Cells.Item(7, 7).Value = "Ají"
$ws.Cells.Item(7, 8).Value = "Inferno"
$ws.Cells.Item(7, 9).Value = "Primera"
$ws.Cells.Item(7, 10).Value = 45
$ws.Cells.Item(7, 11).Value = 16000
$ws.Cells.Item(7, 12).Value = 17000
$ws.Cells.Item(7, 13).Value = 16667
$ws.Cells.Item(7, 14).Value = "`$/caja 12 kilos"
$ws.Cells.Item(7, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(7, 16).Value = 1389
$ws.Cells.Item(7, 17).Value = 12
$ws.Cells.Item(7, 18).Value = "Hortaliza"
